# 023FW_DataFlow/01.introduction.pptx : "update 023FW_DataFlow form 01 to 03."
#
# The subtitle on slide 1 reads "ZOMI酱" (two runs: "ZOMI" in one run and
# "酱" in a second, differently-tagged run). The commit drops the trailing
# "酱" run so the subtitle just reads "ZOMI".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape 2 ("副标题 2") holds the "ZOMI酱" subtitle text.
$sh = $s.Shapes.Item(2)
$sh.TextFrame.TextRange.Text = "ZOMI"
